# "Fixed 'Definition' typo" — also re-homes the "Ban a User" glossary
# entry next to the other Concept rows and fixes its own
# "operation" -> "Operation" capitalization typo along the way
# (matches the author's resorted Glossary sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the column header typo: "Defination" -> "Definition"
$ws.Range("C1").Value = "Definition"

# 2. Relocate the "Ban a User" row (currently row 55, with a blank
#    separator row above it at 54) so it sits right after the
#    "Concept" rows, keeping a blank separator row below it too.
#    Insert two fresh rows at 13 (old row 13 "uid" onward shifts down
#    by two, so the row + its old separator land at 56/57).
$ws.Rows("13:14").Insert()

# Copy the relocated row (now at 57) into the new row 13; row 14
# stays blank, mirroring the original section-separator gap.
$ws.Rows("57:57").Copy($ws.Rows("13:13"))

# Remove the vacated row and its now-redundant blank separator.
$ws.Rows("56:57").Delete()

# 3. Fix the capitalization typo in the moved row's Category cell.
$ws.Range("B13").Value = "Operation"
